$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.108.26"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "3.770.13"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "628.01"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.96"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("D7").Value = "3.767.98"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").Value = "  +1.92%  "
$ws.Range("E12").Value = "  -3.46%  "
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.90"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "4.407.50"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "3.760.35"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").Value = "69.118.43"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.56"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "463.18"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.52"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.96"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.98"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "3.919.82"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.67"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("E33").Value = "  -2.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.52"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.170"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.34%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.99"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("B38").Value = "RenzoRestakedETH"
$ws.Range("C38").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D38").Value = "3.725.06"
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.31"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.957"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "157.66"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.43"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("E47").Value = "  +3.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.12"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.65"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("E51").Value = "  -0.49%  "
